$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels in row 2 (new "IMF (20%)" columns, shifted "IMF" columns) ---
$ws.Range("D2").Value = "IMF (20%) - Sales"
$ws.Range("E2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("F2").Value = "IMF - Sales"
$ws.Range("G2").Value = "IMF - Sales + Emp"
$ws.Range("L2").Value = "IMF (20%) - Sales"
$ws.Range("M2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("N2").Value = "IMF - Sales"
$ws.Range("O2").Value = "IMF - Sales + Emp"
$ws.Range("T2").Value = "IMF (20%) - Sales"
$ws.Range("U2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("V2").Value = "IMF - Sales"
$ws.Range("W2").Value = "IMF - Sales + Emp"

# --- Update data values for rows 4, 6, 7, 8 ---
# Row 4
$ws.Range("D4").Value = 1.090355678702299
$ws.Range("E4").Value = 0.9036458532427226
$ws.Range("F4").Value = 5.451778393511495
$ws.Range("G4").Value = 4.518229266213623
$ws.Range("L4").Value = 0.2133331070738557
$ws.Range("M4").Value = 0.2133203841653253
$ws.Range("N4").Value = 0.2133331070738557
$ws.Range("O4").Value = 0.2133203841653253
$ws.Range("T4").Value = 956175459402
$ws.Range("U4").Value = 956426351018
$ws.Range("V4").Value = 956175459402
$ws.Range("W4").Value = 956426351018

# Row 6
$ws.Range("D6").Value = 0.7050127574638287
$ws.Range("E6").Value = 2.773577005388292
$ws.Range("F6").Value = 3.525063787319145
$ws.Range("G6").Value = 13.86788502694147
$ws.Range("L6").Value = 1.072913968120253
$ws.Range("M6").Value = 0.4097340342028742
$ws.Range("N6").Value = 1.072913968120253
$ws.Range("O6").Value = 0.4097340342028742
$ws.Range("T6").Value = 2124037026
$ws.Range("U6").Value = 17056843395
$ws.Range("V6").Value = 2124037026
$ws.Range("W6").Value = 17056843395

# Row 7
$ws.Range("D7").Value = 0.6942562708641375
$ws.Range("E7").Value = 0.7822019816789091
$ws.Range("F7").Value = 3.471281354320685
$ws.Range("G7").Value = 3.911009908394546
$ws.Range("L7").Value = 0.2339064661756046
$ws.Range("M7").Value = 0.3977401926186418
$ws.Range("N7").Value = 0.2339064661756046
$ws.Range("O7").Value = 0.3977401926186418
$ws.Range("T7").Value = 1750040641
$ws.Range("U7").Value = 868911225
$ws.Range("V7").Value = 1750040641
$ws.Range("W7").Value = 868911225

# Row 8
$ws.Range("D8").Value = 0.2250293725144459
$ws.Range("E8").Value = 0.5258857082920571
$ws.Range("F8").Value = 1.125146862572229
$ws.Range("G8").Value = 2.629428541460285
$ws.Range("L8").Value = 0.4134566240721217
$ws.Range("M8").Value = 0.3808706636599813
$ws.Range("N8").Value = 0.4134566240721217
$ws.Range("O8").Value = 0.3808706636599813
$ws.Range("T8").Value = 41286747272
$ws.Range("U8").Value = 47534570308
$ws.Range("V8").Value = 41286747272
$ws.Range("W8").Value = 47534570308

